$d = $word.ActiveDocument

$replacements = @(
    @("252×3=", "552×7="),
    @("820×6=", "259×3="),
    @("888×6=", "140×8="),
    @("230×4=", "438×2="),
    @("933×5=", "562×9="),
    @("979×6=", "678×4="),
    @("540×9=", "312×6="),
    @("416×8=", "842×4="),
    @("199×3=", "104×9="),
    @("734×5=", "328×5="),
    @("994×2=", "528×9="),
    @("566×7=", "627×6="),
    @("499×2=", "520×4="),
    @("610×4=", "874×6="),
    @("463×8=", "189×6="),
    @("359×9=", "890×5="),
    @("979×5=", "332×3="),
    @("896×8=", "411×6="),
    @("927×5=", "775×2="),
    @("638×7=", "631×6="),
    @("561×7=", "415×2="),
    @("578×4=", "616×3="),
    @("444×5=", "783×7="),
    @("674×5=", "416×6="),
    @("937×3=", "850×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
